$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reassign the Experience/Time/Clicks "blocks" among participants 3-6 ---
# new participant 3 <- old participant 4 block
$ws.Range("B10").Value = "None"
$ws.Range("D10").Value = 15.1
$ws.Range("E10").Value = 4

$ws.Range("B11").Value = "None"
$ws.Range("D11").Value = 11.13
$ws.Range("E11").Value = 10

$ws.Range("B12").Value = "None"
$ws.Range("D12").Value = 10.61
$ws.Range("E12").Value = 4

$ws.Range("B13").Value = "None"
$ws.Range("D13").Value = 15.39
$ws.Range("E13").Value = 2

# new participant 4 <- old participant 6 block
$ws.Range("B14").Value = "None"
$ws.Range("D14").Value = 31.39
$ws.Range("E14").Value = 23

$ws.Range("B15").Value = "None"
$ws.Range("D15").Value = 35.88
$ws.Range("E15").Value = 17

$ws.Range("B16").Value = "None"
$ws.Range("D16").Value = 9.1
$ws.Range("E16").Value = 4

$ws.Range("B17").Value = "None"
$ws.Range("D17").Value = 7.63
$ws.Range("E17").Value = 4

# new participant 5 <- old participant 3 block
$ws.Range("B18").Value = "Some"
$ws.Range("D18").Value = 11.2
$ws.Range("E18").Value = 7

$ws.Range("B19").Value = "Some"
$ws.Range("D19").Value = 15.51
$ws.Range("E19").Value = 8

$ws.Range("B20").Value = "Some"
$ws.Range("D20").Value = 11.27
$ws.Range("E20").Value = 5

$ws.Range("B21").Value = "Some"
$ws.Range("D21").Value = 10.4
$ws.Range("E21").Value = 2

# new participant 6 <- old participant 5 block
$ws.Range("B22").Value = "Some"
$ws.Range("D22").Value = 7.88
$ws.Range("E22").Value = 4

$ws.Range("B23").Value = "Some"
$ws.Range("D23").Value = 11.91
$ws.Range("E23").Value = 8

$ws.Range("B24").Value = "Some"
$ws.Range("D24").Value = 8.15
$ws.Range("E24").Value = 3

$ws.Range("B25").Value = "Some"
$ws.Range("D25").Value = 5.87
$ws.Range("E25").Value = 2

# --- Clear the (redundant) explicit fill/border style on participant 7 & 8's rows 2-4 ---
$ws.Range("B27:E27").ClearFormats()
$ws.Range("B28:E28").ClearFormats()
$ws.Range("B29:E29").ClearFormats()
$ws.Range("B31:E31").ClearFormats()
$ws.Range("B32:E32").ClearFormats()
$ws.Range("B33:E33").ClearFormats()

# --- Update the active selection ---
$ws.Range("G37").Select()
